# Apply the "add rows" diff to the coffee-log workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Minor style tweak already present on row 4: G4 switches from the
#    "no reading-order" general body style (same as C4..N4) to the
#    "reading-order" general body style already used by B4.
# ------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("G4").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Build the brand-new "pink/bold/Comic Sans MS" font+style used by
#    the whole new row 5. We derive it from the existing body styles
#    (fontId=2, Comic Sans MS / theme color) so that only Bold + Color
#    need to change - this keeps the new font definition identical to
#    <b/><color rgb="FFFF00FF"/><name val="Comic Sans MS"/>.
# ------------------------------------------------------------------

# General (non-date) seed cell -> becomes the template for most of row 5
$ws.Range("P1").Value = "seed"
$ws.Range("C4").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Font.Bold = $true
$ws.Range("P1").Font.Color = 16711935

# Date seed cell -> becomes the template for A5 (keeps yyyy/m/d format)
$ws.Range("P2").Value = 1
$ws.Range("A4").Copy()
$ws.Range("P2").PasteSpecial(-4122)
$ws.Range("P2").Font.Bold = $true
$ws.Range("P2").Font.Color = 16711935

# ------------------------------------------------------------------
# 3) Write row 5 values.
# ------------------------------------------------------------------
$ws.Range("A5").Value = 44962.0
$ws.Range("B5").Value = "Colonbia Willa"
$ws.Range("C5").Value = "Dark"
$ws.Range("D5").Value = "弱い"
$ws.Range("E5").Value = "ふつう"
$ws.Range("F5").Value = "ユニーク"
$ws.Range("G5").Value = 12.0
$ws.Range("H5").Value = 10.0
$ws.Range("I5").Value = 100.0
$ws.Range("J5").Value = 92.0
$ws.Range("K5").Value = "よつ葉牛乳"
$ws.Range("L5").Value = 100.0
$ws.Range("M5").Value = "キビ砂糖"
$ws.Range("N5").Value = 5.0

# ------------------------------------------------------------------
# 4) Apply the new pink/bold/Comic Sans MS format across row 5.
# ------------------------------------------------------------------
$ws.Range("P1").Copy()
$ws.Range("B5:N5").PasteSpecial(-4122)

$ws.Range("P2").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Clean up the helper/scratch cells used to build the new style so no
# stray formatted-but-empty cells remain outside the used range.
$ws.Range("P1:P2").Clear()

# ------------------------------------------------------------------
# 5) Column width adjustments (auto-fit side effects of the new row).
#    ColumnWidth (character units) is offset from the stored OOXML
#    width by 5/6, so we compensate to land as close as possible to
#    the target widths of 12.38 / 11.63 / 19.38.
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11.5
$ws.Columns.Item(7).ColumnWidth = 10.833333333333332
$ws.Columns.Item(13).ColumnWidth = 18.5
